# "Add files via upload" — refresh the VLOOKUP source table on Sheet3 with a
# new day's numbers (03-nov), which ripples into the existing CB/CC VLOOKUP
# columns on Sheet1, and append a new "03-nov" column (CH) on Sheet1 that
# snapshots that day's looked-up values.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws3 = $wb.Worksheets.Item("Sheet3")

# 1) Update the lookup table (Sheet3 A20:B36) with the new day's values.
#    This table feeds Sheet3!C2:C18 (IFERROR/VLOOKUP) which in turn feeds
#    Sheet1's CB/CC VLOOKUP columns - everything recalculates from here.
$newValues = @{
    20 = 12.529154809458598    # 3D QUESO 92GX27
    21 = 0.9860418682175317    # CHEETOS 94GRX24
    22 = 8.6477451095384463    # DORITOS QUESO 140GX19
    23 = 5.8590180893980532    # DORITOS QUESO 70X40G
    24 = 9.2131255648029455    # DORITOS QUESO 77GX26
    25 = 0                     # LAYS CEBOLLA CARAMELIZADA 85GX25
    26 = 3.9168015747708331    # LAYS CLASICAS 145GRX18
    27 = 12.713886190011543    # LAYS CLASICAS 249GRX14
    28 = 6.9465275606272439    # LAYS CLASICAS 40GX68
    29 = 4.8284998318561394    # LAYS CLASICAS 94GRX25
    30 = 0                     # LAYS ONDAS FH 30GX72
    31 = 2.4464693019445143    # LAYS ONDAS FH 70GX28
    32 = 8.894002293627068     # LAYS QSO Y CEBOLLA 34GX72
    33 = 3.3420105526665189    # PEHUAMAR ACANALADA 520GX9
    34 = 5.5218347709050883    # PEHUAMAR MAICITOS 285GX10
    35 = 8.813784465112489     # PEHUAMAR PAPA LISA 520GX9
    36 = 42.320547261257438    # QUAKER AVENA INSTANT FORTIF 18X280G
}

foreach ($row in $newValues.Keys) {
    $ws3.Cells.Item($row, 2).Value = $newValues[$row]
}

# 2) Add the new "03-nov" column (CH) on Sheet1: header in row 1, and for
#    each data row (2-18) the value now produced by the CB/CC VLOOKUP for
#    that row's product (same number that Sheet3!C<row> now computes).
#    Match the formatting of the existing last column (CG) so the new
#    column looks like the others.
$ws1.Range("CH1").Value = "03-nov"
$ws1.Range("CH1").NumberFormat = $ws1.Range("CG1").NumberFormat()

for ($r = 2; $r -le 18; $r++) {
    $ws1.Cells.Item($r, 86).Value = $ws1.Cells.Item($r, 80).Value()
    $ws1.Cells.Item($r, 86).NumberFormat = $ws1.Cells.Item($r, 85).NumberFormat()
}

# 3) Restore the prior "last touched" selection state recorded in the diff.
$ws1.Range("CH22").Select() | Out-Null
